$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("A1").Value = "Function Category"

# Add new column L header/value
$ws.Range("L1").Value = "ETC"
$ws.Range("L2").Value = "X"

# Update row 2 values (A2, B2) to new Korean labels
$ws.Range("A2").Value = "기능 이름"
$ws.Range("B2").Value = "요약"
